$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.339.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.797.85'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +21.98%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '618.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.91'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.795.05'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +21.89%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.555'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.35%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +11.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.39'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.84'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +12.07%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.431.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +22.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.787.20'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +21.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.589.09'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.89%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '526.67'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.89'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +23.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.749'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.63'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.91%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +10.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.59'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.66%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +30.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.53'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +9.47%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.30'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +14.82%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.26%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.22'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +11.45%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +10.85%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.95%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.86%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.75'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.07%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '430.94'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +15.54%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.91'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.142.75'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +12.84%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '44.25'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.06%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.96%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.01%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +9.05%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.41'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.01%  '
